$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("web 125")

# Header row additions (new quiz columns)
$ws.Range("G1").Value = "Q02"
$ws.Range("H1").Value = "Q03"
$ws.Range("I1").Value = "Q04"

# Row 2
$ws.Range("G2").Formula = "=(5/15)*10"
$ws.Range("H2").Value = 0
$ws.Range("I2").Formula = "=(10/25)*10"

# Row 3
$ws.Range("G3").Value = 0
$ws.Range("H3").Value = 0
$ws.Range("I3").Value = 0

# Row 4
$ws.Range("G4").Formula = "=(7/15)*10"
$ws.Range("H4").Value = 0
$ws.Range("I4").Formula = "=(12/25)*10"

# Row 5
$ws.Range("G5").Value = 0
$ws.Range("H5").Value = 0
$ws.Range("I5").Value = 0

# Row 6
$ws.Range("G6").Formula = "=(12/15)*10"
$ws.Range("H6").Formula = "=(16.5/30)*10"
$ws.Range("I6").Formula = "=(22/25)*10"

# Row 7
$ws.Range("G7").Value = 0
$ws.Range("H7").Value = 0
$ws.Range("I7").Value = 0

# Row 12
$ws.Range("G12").Value = 0
$ws.Range("H12").Value = 0
$ws.Range("I12").Value = 0

# Row 13
$ws.Range("G13").Value = 0
$ws.Range("H13").Value = 0
$ws.Range("I13").Value = 0

# Row 15
$ws.Range("G15").Value = 0
$ws.Range("H15").Value = 0
$ws.Range("I15").Value = 0

# Row 16
$ws.Range("G16").Formula = "=(9/15)*10"
$ws.Range("H16").Formula = "=(25.5/30)*10"
$ws.Range("I16").Formula = "=(18/25)*10"

# Row 17
$ws.Range("G17").Formula = "=(15/15)*10"
$ws.Range("H17").Formula = "=(25.5/30)*10"
$ws.Range("I17").Formula = "=(23/25)*10"

# Update selection to match final author position
$ws.Range("I31").Select()
